# Generate Report for Handoff
#
# This updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# columns for the row corresponding to file
# "6781aca2-82b6-4d69-94ea-5dfa5d53cef7.md" on all three report sheets,
# reflecting a freshly (re-)generated handoff report.

$wb = $excel.ActiveWorkbook

# --- Sheet "Overview": column G is "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-08-28 16:42:12"

# --- Sheet "zh-cn": column H is "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-08-28 16:42:07"

# --- Sheet "de-de": column H is "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-08-28 16:42:12"
